# DSS.xlsx -- "Add files via upload"
#
# Appends 16 new training-record rows (two new employees, 8 course rows
# each) to the bottom of Sheet1, right after the existing data that ends
# at row 1180. Also widens the Print_Area defined name to cover the new
# rows and leaves the sheet's scroll/selection state the way the author
# left it (selection on B1201).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Stamp the row-level formatting (fill/border/font/number-format) for
#    the new rows by copying it from existing rows that already carry the
#    exact style combination we need, so every new cell ends up with the
#    same style index as the rest of the sheet (no ad-hoc new styles).
# ---------------------------------------------------------------------

# Row 1181-1188 ("Mohamed ...") : col A keeps the "s=10" look, columns
# B:E pick up the "s=11" look used by every other row in that position.
$srcA10 = $ws.Range("A1144")
$srcBCDE11 = $ws.Range("B1152:E1152")

for ($r = 1181; $r -le 1188; $r++) {
    $srcA10.Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $srcBCDE11.Copy()
    $ws.Range("B$r" + ":E$r").PasteSpecial(-4122)
}

# Row 1189 ("Youssef ...", first course row): whole row mirrors 1144
# (A:C s=10, D s=44, E s=17).
$srcRow1144 = $ws.Range("A1144:E1144")
$srcRow1144.Copy()
$ws.Range("A1189:E1189").PasteSpecial(-4122)

# Rows 1190-1196 (remaining "Youssef ..." course rows): whole row mirrors
# 1145 (A:E all s=10/17).
$srcRow1145 = $ws.Range("A1145:E1145")
for ($r = 1190; $r -le 1196; $r++) {
    $srcRow1145.Copy()
    $ws.Range("A$r" + ":E$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Write the cell values/content for the 16 new rows.
# ---------------------------------------------------------------------

$courses = @(
    "30 Hours Construction Safety & Health",
    "30 Hours G. Industry Safety & Health",
    "Electrical Safety & LOTO",
    "Fire Marshal",
    "Scaffold Competent Person",
    "Lifting & Rigging Competent Person",
    "Health & Safety Risk Assessment",
    "Safety Management System & PTW"
)
$dates = @(
    "05-01-2025",
    "10-01-2025",
    "06-01-2025",
    "03-01-2025",
    "01-01-2025",
    "02-01-2025",
    "07-01-2025",
    "08-01-2025"
)

# Block 1: DSS2180-DSS2187, Mohamed Ahmed Mohamed Mahmoud -> rows 1181-1188
$name1 = "Mohamed Ahmed Mohamed Mahmoud"
$badgeStart1 = 2180
for ($i = 0; $i -lt 8; $i++) {
    $r = 1181 + $i
    $ws.Cells.Item($r, 1).Value = "DSS" + [string]($badgeStart1 + $i)
    $ws.Cells.Item($r, 2).Value = $name1
    $ws.Cells.Item($r, 3).Value = $courses[$i]
    $ws.Cells.Item($r, 4).Value = $dates[$i]
    $ws.Cells.Item($r, 5).Value = 1
}

# Block 2: DSS2188-DSS2195, Youssef Ahmed Mohamed Mahmoud -> rows 1189-1196
$name2 = "Youssef Ahmed Mohamed Mahmoud"
$badgeStart2 = 2188
for ($i = 0; $i -lt 8; $i++) {
    $r = 1189 + $i
    $ws.Cells.Item($r, 1).Value = "DSS" + [string]($badgeStart2 + $i)
    $ws.Cells.Item($r, 2).Value = $name2
    $ws.Cells.Item($r, 3).Value = $courses[$i]
    $ws.Cells.Item($r, 4).Value = $dates[$i]
    $ws.Cells.Item($r, 5).Value = 1
}

# ---------------------------------------------------------------------
# 3. Grow the printed area to match the new last row of data.
# ---------------------------------------------------------------------

$ws.PageSetup.PrintArea = '$A$1:$E$1207'

# ---------------------------------------------------------------------
# 4. Leave the sheet scrolled/selected the way it was saved: cursor on
#    B1201.
# ---------------------------------------------------------------------

$ws.Range("B1201").Select()

"Appended rows 1181:1196 and updated Print_Area."
